# Update countries & provincias Spain
# - Swap ranking order of Singapur / Bielorrusia (row 29 / row 30) and
#   refresh several countries' case figures.
# - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 29 / 30: Singapur overtakes Bielorrusia in the ranking ---------
# Row 29 now holds Singapur's (updated) figures, row 30 now holds
# Bielorrusia's (unchanged) figures - i.e. the two rows swap places.
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 26098
$ws.Range("C29").Value = 752
$ws.Range("D29").Value = 4809
$ws.Range("E29").Value = 21268
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 21

$ws.Range("A30").Value = "Bielorrusia"
$ws.Range("B30").Value = 25825
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 7711
$ws.Range("E30").Value = 17968
$ws.Range("F30").Value = 92
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 146

# --- Row 34: updated figures --------------------------------------------
$ws.Range("D34").Value = 6696
$ws.Range("E34").Value = 9647

# --- Row 65: updated figures ----------------------------------------------
$ws.Range("B65").Value = 4341
$ws.Range("C65").Value = 322
$ws.Range("D65").Value = 1303
$ws.Range("E65").Value = 3021

# --- Row 119: updated figures ---------------------------------------------
$ws.Range("B119").Value = 652
$ws.Range("C119").Value = 5
$ws.Range("D119").Value = 383
$ws.Range("E119").Value = 257
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 12

# --- Row 128: updated figures ----------------------------------------------
$ws.Range("D128").Value = 383
$ws.Range("E128").Value = 50

# --- Timestamp banner in A1 ------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 09:35"
